# "Saving account Test cases added"
#
# NewCheckingAccount (sheet1) already contains the checking-account test
# data. This adds an equivalent "NewSavingAccount" sheet (currently the
# blank "Sheet2") populated with saving-account test rows, makes it the
# active/selected tab, and leaves the checking-account sheet selection
# pointed at its data range.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # NewCheckingAccount
$ws2 = $wb.Worksheets.Item(2)   # Sheet2 -> NewSavingAccount

# 1. Rename the second sheet.
$ws2.Name = "NewSavingAccount"

# 2. Bring over the same cell formatting (borders / header fill) used by
#    the checking account sheet so the new sheet looks the same.
$ws1.Range("A1:B3").Copy()
$ws2.Range("A1:B3").PasteSpecial(-4122)   # xlPasteFormats

# 3. Match the column widths used on sheet1.
$ws2.Columns.Item(1).ColumnWidth = 13.166666666666666
$ws2.Columns.Item(2).ColumnWidth = 14.736979166666666

# 4. Fill in the saving-account test data.
$ws2.Range("A1").Value = "Account Name"
$ws2.Range("B1").Value = "Deposit Amount"
$ws2.Range("A2").Value = "CCC"
$ws2.Range("B2").Value = 4000
$ws2.Range("A3").Value = "DDD"
$ws2.Range("B3").Value = 5000

# 5. Update selections: checking account sheet selects its whole table,
#    saving account sheet becomes the active tab with C3 selected.
$ws1.Range("A1:B3").Select()
$ws2.Activate()
$ws2.Range("C3").Select()
